$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the first data row (2025-11-10), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Append four new data rows for 2026-02-04 .. 2026-02-07.
# Column A must stay a plain text string (not an Excel date serial), matching
# the rest of the sheet, so format the range as text before typing the
# values and then clear the formatting again (keeps default style s="0").
$dateRange = $ws.Range("A87:A90")
$dateRange.NumberFormat = "@"

$ws.Range("A87").Value = "2026-02-04"
$ws.Range("B87").Value = 320.0
$ws.Range("C87").Value = 220.0
$ws.Range("D87").Value = 351.0

$ws.Range("A88").Value = "2026-02-05"
$ws.Range("B88").Value = 320.0
$ws.Range("C88").Value = 220.0
$ws.Range("D88").Value = 208.0

$ws.Range("A89").Value = "2026-02-06"
$ws.Range("B89").Value = 320.0
$ws.Range("C89").Value = 220.0
$ws.Range("D89").Value = 75.0

$ws.Range("A90").Value = "2026-02-07"
$ws.Range("B90").Value = 320.0
$ws.Range("C90").Value = 220.0
$ws.Range("D90").Value = 329.0

$dateRange.ClearFormats()
